# VOLHO6.xlsx update: "All Datatables updated. MTR366 added."
#
# - EffectiveDate / PreviousExpDate values of "07252023" become "08252023"
#   (shared string used by F2, I2, F3 and I3 on the ho6customerInfo sheet).
# - The sheet's active cell selection moves from F8 to E11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "08252023"
$ws.Range("I2").Value = "08252023"
$ws.Range("F3").Value = "08252023"
$ws.Range("I3").Value = "08252023"

[void]$ws.Range("E11").Select()
